$d = $word.ActiveDocument

# Each pair below is unique in the source document, so a single
# Find/Replace (wdReplaceAll) per pair safely targets exactly one cell.
$pairs = @(
    @("77÷5=", "85÷3="),
    @("39÷9=", "28÷9="),
    @("44÷9=", "83÷8="),
    @("84÷9=", "13÷8="),
    @("30÷2=", "95÷7="),
    @("12÷7=", "21÷4="),
    @("45÷6=", "85÷3="),
    @("98÷8=", "91÷4="),
    @("48÷3=", "65÷7="),
    @("29÷3=", "33÷3="),
    @("14÷7=", "56÷6="),
    @("60÷7=", "46÷6="),
    @("11÷2=", "79÷6="),
    @("88÷9=", "97÷2="),
    @("68÷5=", "82÷4="),
    @("86÷5=", "96÷2="),
    @("65÷9=", "96÷3="),
    @("27÷8=", "89÷3="),
    @("28÷7=", "30÷5="),
    @("99÷4=", "73÷7="),
    @("75÷2=", "36÷2="),
    @("57÷8=", "51÷9="),
    @("10÷9=", "87÷3="),
    @("77÷8=", "43÷7="),
    @("97÷4=", "39÷7="),
)

foreach ($pair in $pairs) {
    $find = $pair[0]
    $replace = $pair[1]
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                            $true, 1, $false, $replace, 2)
}

